# Add the new "Base registration rate" variable to the FTT_variables masterfile:
#  - FTT-Tr sheet: insert a new variable row (row 4)
#  - Time_Horizons sheet: insert the matching time-horizon lookup row (row 10)

$wb = $excel.ActiveWorkbook

# --- FTT-Tr sheet -----------------------------------------------------
$wsTr = $wb.Worksheets.Item("FTT-Tr")

# Insert a new row above the current row 4 (TTRT / road tax), shifting
# everything below down by one.
$wsTr.Rows.Item(4).Insert()

$wsTr.Range("A4").Value = "Base registration rate"
$wsTr.Range("B4").Value = 1
$wsTr.Range("C4").Value = -99
$wsTr.Range("D4").Value = "FTT-Tr registration rate as share of purchase price"
$wsTr.Range("E4").Value = "VTTI"
$wsTr.Range("F4").Value = "TIME"
$wsTr.Range("G4").Value = "RSHORTTI"
$wsTr.Range("H4").Value = 0
$wsTr.Range("I4").Value = "All"

# --- Time_Horizons sheet ------------------------------------------------
$wsTh = $wb.Worksheets.Item("Time_Horizons")

# Insert a new row above the current row 10 (TTRT), shifting everything
# below down by one, and add the time-horizon entry for the new variable.
$wsTh.Rows.Item(10).Insert()

$wsTh.Range("A10").Value = "Base registration rate"
$wsTh.Range("B10").Value = "tl_2001"

# --- Selection / active sheet bookkeeping -------------------------------
# Record the selection on Time_Horizons first, then leave FTT-Tr as the
# active/selected sheet (matching where the edit actually took place).
$wsTh.Activate()
$wsTh.Range("B10").Select()

$wsTr.Activate()
$wsTr.Range("A4").Select()
